$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume figures (D = Price, E = Volume(1h)).
# Price strings that parse as plain numbers are apostrophe-prefixed so Excel
# keeps them as text (matching the source data) instead of coercing to a number.

# Row 2
$ws.Range("D2").Value = '29.221.10'
$ws.Range("E2").Value = '  -0.48%  '

# Row 3
$ws.Range("D3").Value = '1.829.17'
$ws.Range("E3").Value = '  -0.75%  '

# Row 4
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").Value = '''237.83'
$ws.Range("E5").Value = '  -1.01%  '

# Row 6
$ws.Range("D6").Value = '''0.6097'
$ws.Range("E6").Value = '  -3.25%  '

# Row 7
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("D8").Value = '''0.07110'
$ws.Range("E8").Value = '  -4.80%  '

# Row 9
$ws.Range("D9").Value = '''0.2830'
$ws.Range("E9").Value = '  -2.62%  '

# Row 10
$ws.Range("E10").Value = '  -3.89%  '

# Row 11
$ws.Range("D11").Value = '''0.07644'
$ws.Range("E11").Value = '  -1.23%  '

# Row 12
$ws.Range("D12").Value = '1.860.58'
$ws.Range("E12").Value = '  +0.77%  '

# Row 13
$ws.Range("D13").Value = '''4.818'
$ws.Range("E13").Value = '  -3.44%  '

# Row 14
$ws.Range("D14").Value = '''0.6382'
$ws.Range("E14").Value = '  -5.93%  '

# Row 15
$ws.Range("D15").Value = '''0.000009932'
$ws.Range("E15").Value = '  -2.77%  '

# Row 16
$ws.Range("D16").Value = '2.076.37'
$ws.Range("E16").Value = '  -0.96%  '

# Row 17
$ws.Range("D17").Value = '''79.74'
$ws.Range("E17").Value = '  -2.93%  '

# Row 18
$ws.Range("D18").Value = '''5.989'
$ws.Range("E18").Value = '  -4.80%  '

# Row 19
$ws.Range("D19").Value = '29.235.01'
$ws.Range("E19").Value = '  -0.58%  '

# Row 20
$ws.Range("D20").Value = '''230.38'
$ws.Range("E20").Value = '  +0.38%  '

# Row 21
$ws.Range("E21").Value = '  -4.23%  '

# Row 22
$ws.Range("E22").Value = '  +0.13%  '

# Row 23
$ws.Range("D23").Value = '''7.039'
$ws.Range("E23").Value = '  -4.93%  '

# Row 24
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '  +0.22%  '

# Row 25
$ws.Range("D25").Value = '''155.52'
$ws.Range("E25").Value = '  -1.87%  '

# Row 26
$ws.Range("D26").Value = '''8.091'
$ws.Range("E26").Value = '  -4.82%  '

# Row 27
$ws.Range("D27").Value = '''0.1298'
$ws.Range("E27").Value = '  -4.13%  '

# Row 28
$ws.Range("E28").Value = '  -3.84%  '

# Row 29
$ws.Range("D29").Value = '''0.06722'
$ws.Range("E29").Value = '  +2.31%  '

# Row 30
$ws.Range("D30").Value = '''1.489'
$ws.Range("E30").Value = '  +3.46%  '

# Row 31
$ws.Range("D31").Value = '''1.459'
$ws.Range("E31").Value = '  -1.96%  '

# Row 32
$ws.Range("D32").Value = '''3.847'
$ws.Range("E32").Value = '  -5.21%  '

# Row 33
$ws.Range("D33").Value = '''3.823'
$ws.Range("E33").Value = '  -6.10%  '

# Row 34
$ws.Range("D34").Value = '''1.127'
$ws.Range("E34").Value = '  -1.23%  '

# Row 35
$ws.Range("D35").Value = '''1.732'
$ws.Range("E35").Value = '  -5.89%  '

# Row 36
$ws.Range("D36").Value = '''0.6587'
$ws.Range("E36").Value = '  -5.36%  '

# Row 37
$ws.Range("D37").Value = '''2.555'
$ws.Range("E37").Value = '  -0.91%  '

# Row 38
$ws.Range("D38").Value = '1.235.75'
$ws.Range("E38").Value = '  -1.05%  '

# Row 39
$ws.Range("D39").Value = '''2.756'
$ws.Range("E39").Value = '  -2.22%  '

# Row 40
$ws.Range("E40").Value = '  -4.72%  '

# Row 41
$ws.Range("D41").Value = '''6.615'
$ws.Range("E41").Value = '  -2.57%  '

# Row 42
$ws.Range("D42").Value = '''0.9352'
$ws.Range("E42").Value = '  +0.22%  '

# Row 43
$ws.Range("E43").Value = '  +0.14%  '

# Row 44
$ws.Range("D44").Value = '1.988.27'
$ws.Range("E44").Value = '  -1.94%  '

# Row 45
$ws.Range("D45").Value = '''100.84'
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("D46").Value = '''63.67'
$ws.Range("E46").Value = '  -2.90%  '

# Row 47
$ws.Range("E47").Value = '  -2.67%  '

# Row 48
$ws.Range("D48").Value = '''1.635'
$ws.Range("E48").Value = '  -4.93%  '

# Row 49
$ws.Range("D49").Value = '''8.605'
$ws.Range("E49").Value = '  -4.30%  '

# Row 50
$ws.Range("D50").Value = '''0.1087'
$ws.Range("E50").Value = '  -5.15%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '''6.536'
$ws.Range("E51").Value = '  -7.54%  '
